# Add three new columns (Wins, Losses, Ties) holding the team's season record
# alongside the existing per-player roster/statistics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — new column headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from an existing header cell (bold, centered, bordered)
# onto the new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (95 wins, 67 losses, 0 ties) for every player row.
for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = 95
    $ws.Cells.Item($row, 31).Value = 67
    $ws.Cells.Item($row, 32).Value = 0
}
